# Apply corrected model output: flip A3/A9 values on sheets y1 and y2.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("y1")
$ws1.Range("A3").Value = 1
$ws1.Range("A9").Value = 0

$ws2 = $wb.Worksheets.Item("y2")
$ws2.Range("A3").Value = 0
$ws2.Range("A9").Value = 1
